$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (pushes SAN MIGUEL/QUILCATE..SAN PABLO rows down by one)
$ws.Rows("20:20").Insert()

# Fill in the new NIEPOS row (row 20)
$ws.Range("A20").Value = "SAN MIGUEL"
$ws.Range("B20").Value = "NIEPOS"
$ws.Range("C20").Value = -6.92511
$ws.Range("D20").Value = -79.12902
$ws.Range("E20").Value = 0
$ws.Rows("20:20").RowHeight = 13.8

# Build the new bordered/number-formatted style on a scratch cell, then copy
# just the formatting onto C20:D20 (avoids leaving unused intermediate styles
# behind in the style table).
$tmpl = $ws.Range("Z1")
$tmpl.NumberFormat = "#,##0.00"
$tmpl.Borders.LineStyle = 1
$tmpl.Copy()
$ws.Range("C20:D20").PasteSpecial(-4122)
$tmpl.Clear()

$ws.Range("G21").Select()
